# Updated symbol list on Tue Dec 27 13:58:58 UTC 2022 with GitHub Actions
#
# Refreshes the cryptocurrency price/volume snapshot in Sheet1. The "Price"
# column (D) holds numeric-looking values that are stored as plain text in
# the workbook, so we force a text number format before writing them back
# to avoid Excel auto-converting them to floating point numbers. A couple
# of "Volume(1h)" cells (E18, E41) also changed because the "Worst in 24h"
# badge moved from KickToken to One.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2"  "243.43"
Set-TextValue "D3"  "24.06"
Set-TextValue "D4"  "5.377"
Set-TextValue "D5"  "0.05923"
Set-TextValue "D6"  "3.396"
Set-TextValue "D7"  "6.502"
Set-TextValue "D8"  "0.8122"
Set-TextValue "D9"  "0.9401"
Set-TextValue "D10" "0.1422"
Set-TextValue "D11" "0.07405"
Set-TextValue "D12" "0.03090"
Set-TextValue "D13" "0.03056"
Set-TextValue "D14" "0.09335"
Set-TextValue "D15" "3.863"
Set-TextValue "D16" "0.001584"
Set-TextValue "D17" "0.04706"

Set-TextValue "D18" "0.0005981"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue "D19" "0.005936"
Set-TextValue "D22" "0.00007501"
Set-TextValue "D27" "0.0002653"
Set-TextValue "D40" "0.03897"

Set-TextValue "D41" "0.006371"
$ws.Range("E41").Value = "40KickTokenKICK"

Set-TextValue "D42" "0.1072"
Set-TextValue "D43" "0.003201"
Set-TextValue "D44" "0.008869"
Set-TextValue "D45" "0.00005214"
Set-TextValue "D47" "0.6712"
Set-TextValue "D48" "0.002029"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.0002000"
